$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 5
$ws.Range("R2").Value = 2.38
$ws.Range("S2").Value = 1.53
$ws.Range("W2").Value = 13
$ws.Range("AA2").Value = 7
$ws.Range("AB2").Value = 23
$ws.Range("AE2").Value = 23
$ws.Range("AH2").Value = 51
$ws.Range("J3").Value = 1.13
$ws.Range("K3").Value = 6
$ws.Range("N3").Value = 2.5
$ws.Range("O3").Value = 1.5
$ws.Range("R3").Value = 2.1
$ws.Range("S3").Value = 1.67
$ws.Range("U3").Value = 9
$ws.Range("AB3").Value = 17
$ws.Range("R4").Value = 2.1
$ws.Range("S4").Value = 1.67
$ws.Range("AD4").Value = 6
$ws.Range("H5").Value = 2.9
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 1.11
$ws.Range("K5").Value = 6.5
$ws.Range("N5").Value = 2.6
$ws.Range("O5").Value = 1.48
$ws.Range("U5").Value = 8
$ws.Range("W5").Value = 17
$ws.Range("Z5").Value = 6
$ws.Range("AD5").Value = 9
$ws.Range("AE5").Value = 19
$ws.Range("G6").Value = 3.6
$ws.Range("I6").Value = 2
$ws.Range("K6").Value = 8.5
$ws.Range("N6").Value = 2.15
$ws.Range("O6").Value = 1.67
$ws.Range("P6").Value = 1.41
$ws.Range("Q6").Value = 2.62
$ws.Range("R6").Value = 1.91
$ws.Range("S6").Value = 1.8
$ws.Range("X6").Value = 34
$ws.Range("AD6").Value = 6.5
$ws.Range("AE6").Value = 9
$ws.Range("AG6").Value = 17
$ws.Range("AH6").Value = 17
$ws.Range("G8").Value = 3.1
$ws.Range("I8").Value = 3.05
$ws.Range("T8").Value = 5.2
$ws.Range("U8").Value = 13
$ws.Range("V8").Value = 14.5
$ws.Range("W8").Value = 50
$ws.Range("X8").Value = 50
$ws.Range("Y8").Value = 100
$ws.Range("AC8").Value = 300
$ws.Range("AD8").Value = 5.3
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 13.5
$ws.Range("AG8").Value = 45
$ws.Range("AH8").Value = 45
$ws.Range("AI8").Value = 90
$ws.Range("I9").Value = 5.5
$ws.Range("L9").Value = 1.3
$ws.Range("M9").Value = 3.4
$ws.Range("P9").Value = 1.37
$ws.Range("Y9").Value = 29
$ws.Range("Z9").Value = 9.5
$ws.Range("AJ9").Value = 351
$ws.Range("M12").Value = 3.15
$ws.Range("R12").Value = 1.77
$ws.Range("S12").Value = 1.94
$ws.Range("G13").Value = 1.85
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 3.85
$ws.Range("P13").Value = 1.39
$ws.Range("T13").Value = 5.4
$ws.Range("U13").Value = 6.9
$ws.Range("V13").Value = 7.1
$ws.Range("W13").Value = 12.5
$ws.Range("X13").Value = 13
$ws.Range("Y13").Value = 23
$ws.Range("Z13").Value = 8.5
$ws.Range("AA13").Value = 5.6
$ws.Range("AB13").Value = 13
$ws.Range("AC13").Value = 60
$ws.Range("AE13").Value = 17
$ws.Range("AF13").Value = 11
$ws.Range("AI13").Value = 35
$ws.Range("AJ13").Value = 450
$ws.Range("G15").Value = 1.8
$ws.Range("H15").Value = 4.33
$ws.Range("I15").Value = 3.7
$ws.Range("J15").Value = 1.02
$ws.Range("K15").Value = 21
$ws.Range("N15").Value = 1.4
$ws.Range("O15").Value = 2.88
$ws.Range("R15").Value = 1.44
$ws.Range("S15").Value = 2.63
$ws.Range("T15").Value = 12
$ws.Range("U15").Value = 12
$ws.Range("W15").Value = 17
$ws.Range("X15").Value = 12
$ws.Range("Z15").Value = 23
$ws.Range("AA15").Value = 9
$ws.Range("AC15").Value = 34
$ws.Range("AD15").Value = 17
$ws.Range("AJ15").Value = 101
$ws.Range("P16").Value = 1.3
$ws.Range("Q16").Value = 3.4
$ws.Range("R16").Value = 1.7
$ws.Range("S16").Value = 2.05
$ws.Range("T16").Value = 8.5
$ws.Range("Z16").Value = 15
$ws.Range("AI16").Value = 41
$ws.Range("AJ16").Value = 201
$ws.Range("G17").Value = 1.4
$ws.Range("H17").Value = 4.75
$ws.Range("N17").Value = 1.5
$ws.Range("O17").Value = 2.5
$ws.Range("AA17").Value = 9.5
$ws.Range("AD17").Value = 21
$ws.Range("G20").Value = 2.88
$ws.Range("I20").Value = 2.2
$ws.Range("L20").Value = 1.29
$ws.Range("M20").Value = 3.5
$ws.Range("N20").Value = 1.93
$ws.Range("O20").Value = 1.93
$ws.Range("P20").Value = 1.4
$ws.Range("Q20").Value = 2.75
$ws.Range("R20").Value = 1.73
$ws.Range("S20").Value = 2
$ws.Range("W20").Value = 34
$ws.Range("AD20").Value = 8.5
$ws.Range("AE20").Value = 11
$ws.Range("G21").Value = 2.7
$ws.Range("H21").Value = 3.3
$ws.Range("I21").Value = 2.25
$ws.Range("R21").Value = 1.53
$ws.Range("S21").Value = 2.38
$ws.Range("W21").Value = 29
$ws.Range("X21").Value = 21
$ws.Range("AD21").Value = 11
$ws.Range("AE21").Value = 13
$ws.Range("AF21").Value = 9.5
$ws.Range("AG21").Value = 23
$ws.Range("AI21").Value = 23
$ws.Range("G23").Value = 1.48
$ws.Range("H23").Value = 4.75
$ws.Range("J23").Value = 26
$ws.Range("L23").Value = 1.1
$ws.Range("M23").Value = 6.5
$ws.Range("N23").Value = 1.36
$ws.Range("O23").Value = 3
$ws.Range("Z23").Value = 26
$ws.Range("AA23").Value = 11
$ws.Range("AB23").Value = 15
$ws.Range("AD23").Value = 23
$ws.Range("G24").Value = 2.05
$ws.Range("H24").Value = 3.4
$ws.Range("I24").Value = 3.1
$ws.Range("J24").Value = 1.05
$ws.Range("K24").Value = 8.5
$ws.Range("L24").Value = 1.29
$ws.Range("M24").Value = 3.5
$ws.Range("N24").Value = 1.95
$ws.Range("O24").Value = 1.85
$ws.Range("P24").Value = 1.4
$ws.Range("Q24").Value = 2.75
$ws.Range("R24").Value = 1.8
$ws.Range("S24").Value = 1.91
$ws.Range("U24").Value = 10
$ws.Range("V24").Value = 9.5
$ws.Range("W24").Value = 19
$ws.Range("X24").Value = 17
$ws.Range("Z24").Value = 10
$ws.Range("AA24").Value = 7
$ws.Range("AB24").Value = 15
$ws.Range("AD24").Value = 10
$ws.Range("AE24").Value = 17
$ws.Range("AF24").Value = 12
$ws.Range("AG24").Value = 34
$ws.Range("AH24").Value = 26
$ws.Range("AJ24").Value = 600
$ws.Range("G27").Value = 1.98
$ws.Range("H27").Value = 3.1
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = 6.7
$ws.Range("L27").Value = 1.32
$ws.Range("M27").Value = 3.1
$ws.Range("N27").Value = 1.95
$ws.Range("O27").Value = 1.75
$ws.Range("P27").Value = 1.4
$ws.Range("Q27").Value = 2.75
$ws.Range("R27").Value = 1.75
$ws.Range("S27").Value = 1.98
$ws.Range("T27").Value = 7.2
$ws.Range("U27").Value = 9.75
$ws.Range("W27").Value = 18.5
$ws.Range("X27").Value = 15.5
$ws.Range("Y27").Value = 25
$ws.Range("Z27").Value = 6.7
$ws.Range("AA27").Value = 6
$ws.Range("AB27").Value = 13.5
$ws.Range("AC27").Value = 60
$ws.Range("AD27").Value = 10.5
$ws.Range("AE27").Value = 22
$ws.Range("AH27").Value = 40
$ws.Range("AI27").Value = 45
$ws.Range("H29").Value = 3.45
$ws.Range("I29").Value = 4.6
$ws.Range("L29").Value = 1.33
$ws.Range("M29").Value = 2.8
$ws.Range("R29").Value = 1.88
$ws.Range("T29").Value = 6.1
$ws.Range("U29").Value = 7.4
$ws.Range("W29").Value = 13
$ws.Range("Z29").Value = 8.75
$ws.Range("AD29").Value = 11.75
$ws.Range("AF29").Value = 15
$ws.Range("AJ29").Value = 800
